# Commit: "U testu create conformation dinamicka putanja umesto hardkodovane"
# (In the "create confirmation" test, use a dynamic path instead of a hardcoded one.)
#
# The hardcoded "pdf_download_path" column (header in Z1, value in Z2 for the
# first data row) is no longer needed in the test fixture, so both cells are
# cleared. The view's window size/selection are also refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hardcoded "pdf_download_path" header and its hardcoded value
# ("C:\Users\Jelena Bulajic\Downloads") from the fixture.
$ws.Range("Z1").ClearContents()
$ws.Range("Z2").Clear()

# Window was resized/scrolled by the author while editing.
$win = $excel.ActiveWindow
$win.Width = 22188
$win.Height = 9000

# Selection moved as part of the edit.
$ws.Range("Y6").Select()
